$wb = $excel.ActiveWorkbook
$pre = $wb.Worksheets.Item("Pre")

# Duplicate the "Pre" sheet (placed right after it) so the new sheet starts
# out with the exact same headers / number-format styles already wired up,
# then turn it into the "Main" sheet with the fuller WWDT / CTL dataset.
$pre.Copy([System.Reflection.Missing]::Value, $pre)
$ws = $wb.ActiveSheet
$ws.Name = "Main"

# --- Data ------------------------------------------------------------
# col B = sample size (TBT), col C = EY fraction
$rows = @(
    @(2,  263, 0.213),
    @(3,  275, 0.218),
    @(4,  310, 0.222),
    @(5,  315, 0.225),
    @(6,  360, 0.23),
    @(7,  372, 0.232),
    @(8,  187, 0.2065),
    @(9,  190, 0.205),
    @(10, 268, 0.215),
    @(11, 273, 0.213),
    @(12, 320, 0.22699999999999998),
    @(13, 328, 0.225)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
$ws.Range("C2:C13").NumberFormat = "0.00%"

# Row labels
$ws.Range("A2").Value = "WWDT"
$ws.Range("A8").Value = "CTL"

# Re-apply the label font across the whole label column so every row in the
# merged groups shares formatting.
$ws.Range("A2:A13").Font.Name = "Arial"
$ws.Range("A2:A13").Font.Size = 10
$ws.Range("A2:A13").Font.Color = 0

# Group the repeated labels
$ws.Range("A2:A7").Merge() | Out-Null
$ws.Range("A8:A13").Merge() | Out-Null

$ws.Range("D19").Select() | Out-Null
